$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.065.22'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.852.49'
$ws.Range('E3').Value = '  +2.64%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = '''232.81'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '''0.619'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').Value = '''40.77'
$ws.Range('E8').Value = '  +4.49%  '
$ws.Range('D9').Value = '''0.332'
$ws.Range('E9').Value = '  +3.33%  '
$ws.Range('D10').Value = '''0.0693'
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '2.125.37'
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('D13').Value = '''11.45'
$ws.Range('E13').Value = '  +5.70%  '
$ws.Range('D14').Value = '1.855.11'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').Value = '''0.676'
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('D16').Value = '''4.67'
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '35.137.60'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '''70.22'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '0.0₃0792'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').Value = '''240.94'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').Value = '''12.27'
$ws.Range('E21').Value = '  +4.81%  '
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  +2.08%  '
$ws.Range('D25').Value = '''173.34'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('D26').Value = '''7.86'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').Value = '''17.58'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('E28').Value = '  +4.37%  '
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').Value = '''3.96'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = '''3.98'
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').Value = '''1.60'
$ws.Range('E34').Value = '  +23.42%  '
$ws.Range('E35').Value = '  +12.56%  '
$ws.Range('D36').Value = '''0.763'
$ws.Range('E36').Value = '  +13.05%  '
$ws.Range('E37').Value = '  +7.67%  '
$ws.Range('D38').Value = '''1.08'
$ws.Range('E38').Value = '  +13.08%  '
$ws.Range('D39').Value = '''90.41'
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = '1.353.63'
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('E41').Value = '  +3.23%  '
$ws.Range('D42').Value = '''14.68'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').Value = '''2.27'
$ws.Range('E43').Value = '  +3.69%  '
$ws.Range('E44').Value = '  -1.44%  '
$ws.Range('E45').Value = '  +1.94%  '
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').Value = '''6.35'
$ws.Range('E47').Value = '  +3.64%  '
$ws.Range('D48').Value = '2.043.38'
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('E49').Value = '  +19.89%  '
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').Value = '''0.0668'
$ws.Range('E51').Value = '  +0.14%  '
